$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns remain plain text (matches original inlineStr cells)
# so numeric-looking strings (e.g. "314.39") are not auto-converted to numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '28.406.25'
$ws.Range('E2').Value = '  +4.10%  '

$ws.Range('D3').Value = '1.796.34'
$ws.Range('E3').Value = '  +1.13%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').Value = '314.39'
$ws.Range('E5').Value = '  +0.32%  '

$ws.Range('E6').Value = '  -0.03%  '

$ws.Range('D7').Value = '0.5453'
$ws.Range('E7').Value = '  +5.15%  '

$ws.Range('D8').Value = '0.3833'
$ws.Range('E8').Value = '  +3.93%  '

$ws.Range('D9').Value = '0.07608'
$ws.Range('E9').Value = '  +2.88%  '

$ws.Range('D10').Value = '42.65'
$ws.Range('E10').Value = '  -0.51%  '

$ws.Range('E11').Value = '  +3.34%  '

$ws.Range('D12').Value = '0.9999'
$ws.Range('E12').Value = '  -0.09%  '

$ws.Range('D13').Value = '21.12'
$ws.Range('E13').Value = '  +3.12%  '

$ws.Range('D14').Value = '6.188'
$ws.Range('E14').Value = '  +1.93%  '

$ws.Range('D15').Value = '7.403'
$ws.Range('E15').Value = '  +6.26%  '

$ws.Range('D16').Value = '1.793.85'
$ws.Range('E16').Value = '  +1.56%  '

$ws.Range('D17').Value = '91.55'
$ws.Range('E17').Value = '  +2.46%  '

$ws.Range('D18').Value = '0.00001072'
$ws.Range('E18').Value = '  +2.44%  '

$ws.Range('D19').Value = '0.06446'
$ws.Range('E19').Value = '  +0.14%  '

$ws.Range('E20').Value = '  -0.03%  '

$ws.Range('D21').Value = '17.31'
$ws.Range('E21').Value = '  +3.20%  '

$ws.Range('E22').Value = '  +2.37%  '

$ws.Range('D23').Value = '28.415.79'
$ws.Range('E23').Value = '  +4.01%  '

$ws.Range('E24').Value = '  +1.89%  '

$ws.Range('D25').Value = '2.129'
$ws.Range('E25').Value = '  +0.39%  '

$ws.Range('D26').Value = '159.58'
$ws.Range('E26').Value = '  +2.92%  '

$ws.Range('D27').Value = '20.72'
$ws.Range('E27').Value = '  +2.45%  '

$ws.Range('D28').Value = '2.396'
$ws.Range('E28').Value = '  +2.79%  '

$ws.Range('D29').Value = '2.001.47'
$ws.Range('E29').Value = '  +1.49%  '

$ws.Range('D30').Value = '123.31'
$ws.Range('E30').Value = '  +1.71%  '

$ws.Range('D31').Value = '1.128'
$ws.Range('E31').Value = '  +5.66%  '

$ws.Range('E32').Value = '  +5.00%  '

$ws.Range('D33').Value = '5.761'
$ws.Range('E33').Value = '  +3.21%  '

$ws.Range('D34').Value = '3.682'
$ws.Range('E34').Value = '  +1.49%  '

$ws.Range('D35').Value = '0.2385'
$ws.Range('E35').Value = '  +17.90%  '

$ws.Range('D36').Value = '0.06790'
$ws.Range('E36').Value = '  +13.38%  '

$ws.Range('D37').Value = '0.02323'
$ws.Range('E37').Value = '  +3.42%  '

$ws.Range('D38').Value = '5.168'
$ws.Range('E38').Value = '  +6.64%  '

$ws.Range('D39').Value = '8.767'
$ws.Range('E39').Value = '  +8.42%  '

$ws.Range('D40').Value = '11.67'
$ws.Range('E40').Value = '  +3.75%  '

$ws.Range('D41').Value = '0.6405'
$ws.Range('E41').Value = '  +3.98%  '

$ws.Range('B42').Value = 'Frax'
$ws.Range('C42').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D42').Value = '0.9993'
$ws.Range('E42').Value = '  -0.09%  '

$ws.Range('B43').Value = 'WEMIXTOKEN'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = '1.409'
$ws.Range('E43').Value = '  -1.69%  '

$ws.Range('D44').Value = '1.161'
$ws.Range('E44').Value = '  +1.89%  '

$ws.Range('D45').Value = '13.65'
$ws.Range('E45').Value = '  +4.27%  '

$ws.Range('D46').Value = '0.5977'
$ws.Range('E46').Value = '  +3.55%  '

$ws.Range('E47').Value = '  +0.93%  '

$ws.Range('D48').Value = '126.50'
$ws.Range('E48').Value = '  +4.20%  '

$ws.Range('D49').Value = '2.002'
$ws.Range('E49').Value = '  +5.94%  '

$ws.Range('D50').Value = '1.150'
$ws.Range('E50').Value = '  +3.12%  '

$ws.Range('D51').Value = '0.06944'
$ws.Range('E51').Value = '  +3.49%  '

# Restore default (unstyled) cell style so only values differ from the original file
$ws.Range("B2:E51").Style = "Normal"
